# Update "latest" optimisation_result.xlsx output (run 190)

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- "Schedule" sheet: rows 2-4 get refreshed pump-run figures ---
$wsSchedule.Range("B2").Value = 46074.22916666666
$wsSchedule.Range("C2").Value = 5
$wsSchedule.Range("D2").Value = 18.9
$wsSchedule.Range("E2").Value = 795.8981452500001
$wsSchedule.Range("F2").Value = 42.11101297619049
$wsSchedule.Range("B3").Value = 46074.75
$wsSchedule.Range("C3").Value = 11
$wsSchedule.Range("D3").Value = 41.58
$wsSchedule.Range("E3").Value = 402.0109762499999
$wsSchedule.Range("F3").Value = 9.668373647186145
$wsSchedule.Range("A4").Value = 46075.04166666666
$wsSchedule.Range("B4").Value = 46075.25
$wsSchedule.Range("C4").Value = 5
$wsSchedule.Range("D4").Value = 18.9
$wsSchedule.Range("E4").Value = 788.95167
$wsSchedule.Range("F4").Value = 41.74347460317461

# --- "Schedule" sheet: brand-new row 5 (extends dimension to A1:F5) ---
$wsSchedule.Range("A5").Value = 46075.70833333334
$wsSchedule.Range("B5").Value = 46076
$wsSchedule.Range("A5:B5").NumberFormat = $wsSchedule.Range("A4:B4").NumberFormat
$wsSchedule.Range("C5").Value = 7
$wsSchedule.Range("D5").Value = 26.46
$wsSchedule.Range("E5").Value = 1081.6806585
$wsSchedule.Range("F5").Value = 40.87984348072562

# --- "Detailed" sheet: refreshed Price/Type/Pump_Status figures ---
$wsDetailed.Range("E12").Value = "ON"
$wsDetailed.Range("B37").Value = 64.38329
$wsDetailed.Range("B38").Value = 105.93651
$wsDetailed.Range("E38").Value = "OFF"
$wsDetailed.Range("B39").Value = 118.32177
$wsDetailed.Range("C39").Value = "historical"
$wsDetailed.Range("B40").Value = 135.70002
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 133.50611
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 106.05354
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 101.07415
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 105.79
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("B48").Value = 107.0409
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 101.76526
$wsDetailed.Range("B50").Value = 96.82904000000001
$wsDetailed.Range("B51").Value = 94.21093
$wsDetailed.Range("B52").Value = 84.79000000000001
$wsDetailed.Range("E52").Value = "ON"
$wsDetailed.Range("B53").Value = 79.95
$wsDetailed.Range("E53").Value = "ON"
$wsDetailed.Range("E54").Value = "ON"
$wsDetailed.Range("E55").Value = "ON"
$wsDetailed.Range("E56").Value = "ON"
$wsDetailed.Range("E57").Value = "ON"
$wsDetailed.Range("E58").Value = "ON"
$wsDetailed.Range("B59").Value = 79.95059999999999
$wsDetailed.Range("E59").Value = "ON"
$wsDetailed.Range("B60").Value = 79.95059999999999
$wsDetailed.Range("E60").Value = "ON"
$wsDetailed.Range("B61").Value = 84.79000000000001
$wsDetailed.Range("B62").Value = 96.94302999999999
$wsDetailed.Range("E62").Value = "OFF"
$wsDetailed.Range("B63").Value = 100.01
$wsDetailed.Range("E63").Value = "OFF"
$wsDetailed.Range("B64").Value = 87.53771
$wsDetailed.Range("E64").Value = "OFF"
$wsDetailed.Range("B65").Value = 79.95027
$wsDetailed.Range("E65").Value = "OFF"
$wsDetailed.Range("B66").Value = 78
$wsDetailed.Range("E66").Value = "OFF"
$wsDetailed.Range("B67").Value = 79.64982999999999
$wsDetailed.Range("E67").Value = "OFF"
$wsDetailed.Range("B68").Value = 90.89
$wsDetailed.Range("E68").Value = "OFF"
$wsDetailed.Range("B69").Value = 84.84577
$wsDetailed.Range("E69").Value = "OFF"
$wsDetailed.Range("B70").Value = 108.01
$wsDetailed.Range("E70").Value = "OFF"
$wsDetailed.Range("B71").Value = 108.01
$wsDetailed.Range("E71").Value = "OFF"
$wsDetailed.Range("B72").Value = 105
$wsDetailed.Range("E72").Value = "OFF"
$wsDetailed.Range("B73").Value = 84.79000000000001
$wsDetailed.Range("E73").Value = "OFF"
$wsDetailed.Range("B74").Value = 84.79000000000001
$wsDetailed.Range("E74").Value = "OFF"
$wsDetailed.Range("B75").Value = 108.01
$wsDetailed.Range("E75").Value = "OFF"
$wsDetailed.Range("B76").Value = 105.00015
$wsDetailed.Range("E76").Value = "OFF"
$wsDetailed.Range("B77").Value = 108.01
$wsDetailed.Range("E77").Value = "OFF"
$wsDetailed.Range("B78").Value = 104.51402
$wsDetailed.Range("E78").Value = "OFF"
$wsDetailed.Range("B79").Value = 108.01
$wsDetailed.Range("E79").Value = "OFF"
$wsDetailed.Range("B80").Value = 105.00015
$wsDetailed.Range("E80").Value = "OFF"
$wsDetailed.Range("B81").Value = 105.00015
$wsDetailed.Range("E81").Value = "OFF"
$wsDetailed.Range("B82").Value = 115.00185
$wsDetailed.Range("E82").Value = "OFF"
$wsDetailed.Range("B83").Value = 84.79000000000001
$wsDetailed.Range("E83").Value = "OFF"
$wsDetailed.Range("B84").Value = 77.99985
$wsDetailed.Range("B85").Value = 78
$wsDetailed.Range("E85").Value = "ON"
$wsDetailed.Range("B86").Value = 78
$wsDetailed.Range("E86").Value = "ON"
$wsDetailed.Range("B87").Value = 93.12528
$wsDetailed.Range("E87").Value = "ON"
$wsDetailed.Range("B88").Value = 85.92464
$wsDetailed.Range("E88").Value = "ON"
$wsDetailed.Range("B89").Value = 81.64194000000001
$wsDetailed.Range("E89").Value = "ON"
$wsDetailed.Range("B90").Value = 90.48074
$wsDetailed.Range("E90").Value = "ON"
$wsDetailed.Range("B91").Value = 84.36094
$wsDetailed.Range("E91").Value = "ON"
$wsDetailed.Range("E92").Value = "ON"
$wsDetailed.Range("B93").Value = 76.1215
$wsDetailed.Range("E93").Value = "ON"
$wsDetailed.Range("B94").Value = 78
$wsDetailed.Range("E94").Value = "ON"
$wsDetailed.Range("B95").Value = 67.2539
$wsDetailed.Range("E95").Value = "ON"
$wsDetailed.Range("B96").Value = 67.11895
$wsDetailed.Range("E96").Value = "ON"
$wsDetailed.Range("B97").Value = 73.38831999999999
$wsDetailed.Range("E97").Value = "ON"
